# Quiz game: add the next question (programming-language quiz row) to the
# existing quiz sheet. The sheet already has two question rows
# (A1:C1 and A2:C2); this appends a third question in the same
# "question, options, answer" layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Which of these is a programming language?"
$ws.Range("B3").Value = "GUI,CGI,PHP,PHD"
$ws.Range("C3").Value = "PHP"
